# Tidsrapportering.xlsx edit script
# Updates the time report: I4 (Kodning, week 18) from 24 to 29 hours,
# which cascades through the weekly/row/grand-total SUM formulas, and
# moves the active selection on the "Marcus" sheet to where editing left
# off (L13).

$wb = $excel.ActiveWorkbook

# --- Worksheet data edit -------------------------------------------------
$ws = $wb.Worksheets.Item("Marcus")

# Increase the "Kodning" hours logged for the G-column week (I4) by 5,
# which cascades through the dependent SUM formulas (K4, I12, C15).
$ws.Range("I4").Value = 29

# Move the active cell/selection as it was left after editing.
$ws.Range("L13").Select()

$wb.Save()
